$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the now-unused "Flagged" rows (originally rows 3,5,7,9,11,13,15,17,19,21,23,25)
# and compact the "Created as Case Accepted" rows up into rows 2-13, advancing each
# month's date by one row and carrying the case_count along with it (row 7's count
# ends up as 30, picked up from the old row 13 "Flagged" count).

$dates  = @(45292, 45323, 45352, 45383, 45413, 45444, 45474, 45505, 45536, 45566, 45597, 45627)
$counts = @(0, 0, 0, 0, 0, 30, 0, 0, 75, 0, 0, 0)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = "Created as Case Accepted"
    $ws.Cells.Item($r, 3).Value = $counts[$i]
}

# Remove the now-extra rows 14:25 that used to hold the "Flagged" entries.
$ws.Range("A14:C25").EntireRow.Delete()

$wb.Save()
